$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Polynomial_regresion" label to the new combined label.
$ws.Range("A4").Value = "Polynomial_regresion/Liner Regresion"

# 2. Add the new "Gardian Boosting Regresion" results row (row 9).
#    Values are entered left-to-right across B:D first, then A last, to
#    match the original authoring order, and use a leading apostrophe so
#    the numeric-looking figures stay plain text (matching every other
#    metric cell in the sheet) rather than being parsed as numbers.
$ws.Range("B9").Value = "'202.8122865941629"
$ws.Range("C9").Value = "'198488.62338179824"
$ws.Range("D9").Value = "'445.5206206022323"
$ws.Range("B9:D9").Style = "Normal"
$ws.Range("A9").Value = "Gardian Boosting Regresion"

# 3. Widen column A so the longer labels fit (matches bestFit width growth).
$ws.Columns.Item(1).ColumnWidth = 32.17

# 4. Move the active selection, as it was left after the edit.
$ws.Range("A11").Select() | Out-Null
